$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the content for the 16/10 daily report block (rows 32-39) ---
# Row 32: Ke hoach (plan) for the day
$ws.Range("C32").Value = "- Công việc 1: biuld màn hình Chi tiết công việc.`n- Công việc 2: Code API lấy data đổ ra màn hình."
# Row 33: Ket qua dat duoc (results)
$ws.Range("C33").Value = "- Công việc 1: hoàn thành`n- Công việc 2: 70%"
# Row 34: Trang thai (status)
$ws.Range("C34").Value = "Chậm tiến độ"
# Row 35: Van de gap phai 1 (issue 1)
$ws.Range("C35").Value = "- Vấn đề 1: hiển thị danh sách Thảo luận chua chính xác về cả bố cục và data."
# Row 36: Van de gap phai 2 (issue 2)
$ws.Range("C36").Value = "- Vấn đề 2: chưa khắc phục được hiện tượng khi chuyển sang activity chi tiết công việc bị tự đông focus vào ô edittext nhập nội dung  thảo luận"
# Row 37: Giai quyet van de 1
$ws.Range("C37").Value = "- Vấn đề 1: đang tìm cách giải quyết ( dự định chuyển sang làm scrollview)"
# Row 38: Giai quyet van de 2
$ws.Range("C38").Value = "- Vấn đề 2: đang tìm hiểu"
# Row 39: Ke hoach ngay mai
$ws.Range("C39").Value = "giải quyế các vấn đề trên + code chức năng màn hình thêm công việc."

# --- Restore/normalize cell formatting that Value-assignment may have disturbed ---
# C32 & C33 use the orange "plan/result" style (same as C5 donor)
$ws.Range("C5").Copy()
$ws.Range("C32").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("C33").PasteSpecial(-4122)

# C34 & C39 use the plain bordered style (same as C7 donor)
$ws.Range("C7").Copy()
$ws.Range("C34").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C39").PasteSpecial(-4122)

# C35, C37, C38 use the quote-prefixed bordered style (same as C8 donor)
$ws.Range("C8").Copy()
$ws.Range("C35").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C37").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C38").PasteSpecial(-4122)

# C36 switches from the quote-prefixed style to the orange wrap-text style (same as C32)
$ws.Range("C32").Copy()
$ws.Range("C36").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Row 36 grows to a double-height row, and gains Wrap Text formatting across
# columns D:F (previously-blank cells picking up a plain wrap-text style).
$ws.Range("D36:F36").WrapText = $true
$ws.Rows.Item(36).RowHeight = 28.5

# --- Update the view state to reflect scrolling down to the newly-edited block ---
$ws.Range("C39").Select()
